# Reorder existing gene-cluster rows and append newly discovered genes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Rv1640c'
$ws.Cells.Item(2, 2).Value = 5
$ws.Cells.Item(2, 3).Value = 'lysX lysS2 lysU mprF Rv1640c MTCY06H11.04c'
$ws.Cells.Item(2, 4).Value = 'FUNCTION: Catalyzes the production of L-lysyl-tRNA(Lys)transfer and the transfer of a lysyl group from L-lysyl-tRNA(Lys) to membrane-bound phosphatidylglycerol (PG), which produces lysylphosphatidylglycerol (LPG), one of the components of the bacterial membrane with a positive net charge. LPG synthesis contributes to the resistance to cationic antimicrobial peptides (CAMPs) and likely protects M.tuberculosis against the CAMPs produced by competiting microorganisms (bacteriocins). In fact, the modification of anionic phosphatidylglycerol with positively charged L-lysine results in repulsion of the peptides. {ECO:0000269|PubMed:19649276}.'
$ws.Cells.Item(2, 5).Value = 37

$ws.Cells.Item(3, 1).Value = 'Rv0020c'
$ws.Cells.Item(3, 2).Value = 5
$ws.Cells.Item(3, 3).Value = 'fhaA TB39.8 Rv0020c'
$ws.Cells.Item(3, 4).Value = 'FUNCTION: Regulates cell growth and peptidoglycan synthesis by binding to MviN. May inhibit the late stages of peptidoglycan synthesis. {ECO:0000269|PubMed:22275220}.'
$ws.Cells.Item(3, 5).Value = 37

$ws.Cells.Item(4, 1).Value = 'Rv3910'
$ws.Cells.Item(4, 2).Value = 5
$ws.Cells.Item(4, 3).Value = 'mviN Rv3910'
$ws.Cells.Item(4, 4).Value = 'FUNCTION: Essential for cell growth and peptidoglycan synthesis. {ECO:0000269|PubMed:22275220}.'
$ws.Cells.Item(4, 5).Value = 37

$ws.Cells.Item(5, 1).Value = 'Rv2752c'
$ws.Cells.Item(5, 2).Value = 5
$ws.Cells.Item(5, 3).Value = 'rnj Rv2752c'
$ws.Cells.Item(5, 4).Value = 'FUNCTION: An RNase that has 5''-3'' exonuclease and possible endonuclease activity. Involved in maturation of rRNA and in some organisms also mRNA maturation and/or decay (By similarity). Has both beta-lactamase and RNase activity, but the physiological relevance of the beta-lactamase activity, i.e. whether it confers antibiotic resistance, has not been shown (PubMed:21568871). {ECO:0000250, ECO:0000269|PubMed:21568871}.'
$ws.Cells.Item(5, 5).Value = 37

$ws.Cells.Item(6, 1).Value = 'Rv2462c'
$ws.Cells.Item(6, 2).Value = 4
$ws.Cells.Item(6, 3).Value = 'tig Rv2462c MTV008.18c'
$ws.Cells.Item(6, 4).Value = 'FUNCTION: Involved in protein export. Acts as a chaperone by maintaining the newly synthesized protein in an open conformation. Functions as a peptidyl-prolyl cis-trans isomerase (By similarity). {ECO:0000250}.'
$ws.Cells.Item(6, 5).Value = 37

$ws.Cells.Item(7, 1).Value = 'Rv3220c'
$ws.Cells.Item(7, 2).Value = 4
$ws.Cells.Item(7, 3).Value = 'pdtaS Rv3220c'
$ws.Cells.Item(7, 4).Value = 'FUNCTION: Member of the two-component regulatory system PdtaR/PdtaS. Autophosphorylates, probably on a histidine residue, and transfers its phosphate group to PdtaR. {ECO:0000269|PubMed:16026786}.'
$ws.Cells.Item(7, 5).Value = 37

$ws.Cells.Item(8, 1).Value = 'Rv0757'
$ws.Cells.Item(8, 2).Value = 4
$ws.Cells.Item(8, 3).Value = 'phoP Rv0757'
$ws.Cells.Item(8, 4).ClearContents()
$ws.Cells.Item(8, 5).Value = 37

$ws.Cells.Item(9, 1).Value = 'Rv3272'
$ws.Cells.Item(9, 2).Value = 3
$ws.Cells.Item(9, 3).Value = 'Rv3272'
$ws.Cells.Item(9, 4).Value = 'FUNCTION: Probably involved in fatty acid metabolism. Binds to fatty acyl-CoAs of varying carbon chain lengths, with the highest binding affinity for palmitoyl-CoA (C16:0). In vitro, alters the cell wall lipid profile and protects mycobacteria from acidic, oxidative and antibiotic stress. May play a significant role in host-pathogen interaction. {ECO:0000269|PubMed:30342240}.'
$ws.Cells.Item(9, 5).Value = 37

$ws.Cells.Item(10, 1).Value = 'Rv1626'
$ws.Cells.Item(10, 2).Value = 3
$ws.Cells.Item(10, 3).Value = 'pdtaR Rv1626'
$ws.Cells.Item(10, 4).Value = 'FUNCTION: Member of the two-component regulatory system PdtaR/PdtaS. {ECO:0000269|PubMed:16026786}.'
$ws.Cells.Item(10, 5).Value = 37

$ws.Cells.Item(11, 1).Value = 'Rv3459c'
$ws.Cells.Item(11, 2).Value = 3
$ws.Cells.Item(11, 3).Value = 'rpsK Rv3459c MTCY13E12.12c'
$ws.Cells.Item(11, 4).Value = 'FUNCTION: Located on the platform of the 30S subunit, it bridges several disparate RNA helices of the 16S rRNA. Forms part of the Shine-Dalgarno cleft in the 70S ribosome. {ECO:0000255|HAMAP-Rule:MF_01310}.'
$ws.Cells.Item(11, 5).Value = 37

$ws.Cells.Item(12, 1).Value = 'Rv3273'
$ws.Cells.Item(12, 2).Value = 3
$ws.Cells.Item(12, 3).Value = 'Rv3273'
$ws.Cells.Item(12, 4).ClearContents()
$ws.Cells.Item(12, 5).Value = 37

$ws.Cells.Item(13, 1).Value = 'Rv2368c'
$ws.Cells.Item(13, 2).Value = 2
$ws.Cells.Item(13, 3).Value = 'Rv2368c MTCY27.12'
$ws.Cells.Item(13, 4).ClearContents()
$ws.Cells.Item(13, 5).Value = 37

$ws.Cells.Item(14, 1).Value = 'Rv0517'
$ws.Cells.Item(14, 2).Value = 2
$ws.Cells.Item(14, 3).Value = 'Rv0517'
$ws.Cells.Item(14, 4).ClearContents()
$ws.Cells.Item(14, 5).Value = 37

$ws.Cells.Item(15, 1).Value = 'Rv0758'
$ws.Cells.Item(15, 2).Value = 2
$ws.Cells.Item(15, 3).Value = 'phoR Rv0758'
$ws.Cells.Item(15, 4).ClearContents()
$ws.Cells.Item(15, 5).Value = 37

$ws.Cells.Item(16, 1).Value = 'Rv3136'
$ws.Cells.Item(16, 2).Value = 1
$ws.Cells.Item(16, 3).Value = 'Rv3136A'
$ws.Cells.Item(16, 4).ClearContents()
$ws.Cells.Item(16, 5).Value = 37

$ws.Cells.Item(17, 1).Value = 'Rv3136'
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(17, 3).Value = 'PPE51 Rv3136'
$ws.Cells.Item(17, 4).ClearContents()
$ws.Cells.Item(17, 5).Value = 37

$ws.Cells.Item(18, 1).Value = 'Rv3135'
$ws.Cells.Item(18, 2).Value = 1
$ws.Cells.Item(18, 3).Value = 'PPE50 Rv3135'
$ws.Cells.Item(18, 4).ClearContents()
$ws.Cells.Item(18, 5).Value = 37

$ws.Cells.Item(19, 1).Value = 'Rv3295'
$ws.Cells.Item(19, 2).Value = 1
$ws.Cells.Item(19, 3).Value = 'Rv3295'
$ws.Cells.Item(19, 4).ClearContents()
$ws.Cells.Item(19, 5).Value = 37
